# Fix Cuisine / Skill / Time values that were shifted one column to the
# right in rows 4 and 6 of the "Recipes" sheet (C = Cuisine, D = Skill,
# E = Time).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recipes")

# Row 4: Quinoa Salad Bowl
$ws.Range("C4").Value = "Intermediate"
$ws.Range("D4").Value = "40 Mins"
$ws.Range("E4").Value = "Western"

# Row 6: Chicken Curry
$ws.Range("C6").Value = "Hard"
$ws.Range("D6").Value = "2 Hours"
$ws.Range("E6").Value = "Indian"
